# Fix Bug in CVRP_L: correct the route/cluster id (col F) and the
# "route start" flag (col E) for rows 3-119 on the active sheet,
# matching the values recorded in the upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 6).Value = 3  # F3: 5 -> 3
$ws.Cells.Item(4, 6).Value = 3  # F4: 5 -> 3
$ws.Cells.Item(5, 5).Value = 1  # E5: 0 -> 1
$ws.Cells.Item(5, 6).Value = 3  # F5: 5 -> 3
$ws.Cells.Item(6, 6).Value = 3  # F6: 5 -> 3
$ws.Cells.Item(7, 6).Value = 3  # F7: 5 -> 3
$ws.Cells.Item(8, 6).Value = 3  # F8: 5 -> 3
$ws.Cells.Item(9, 5).Value = 1  # E9: 0 -> 1
$ws.Cells.Item(9, 6).Value = 3  # F9: 5 -> 3
$ws.Cells.Item(10, 6).Value = 4  # F10: 10 -> 4
$ws.Cells.Item(11, 6).Value = 4  # F11: 10 -> 4
$ws.Cells.Item(12, 5).Value = 1  # E12: 0 -> 1
$ws.Cells.Item(12, 6).Value = 4  # F12: 10 -> 4
$ws.Cells.Item(13, 6).Value = 4  # F13: 10 -> 4
$ws.Cells.Item(14, 6).Value = 4  # F14: 10 -> 4
$ws.Cells.Item(15, 6).Value = 4  # F15: 10 -> 4
$ws.Cells.Item(16, 5).Value = 1  # E16: 0 -> 1
$ws.Cells.Item(16, 6).Value = 4  # F16: 10 -> 4
$ws.Cells.Item(17, 5).Value = 1  # E17: 0 -> 1
$ws.Cells.Item(17, 6).Value = 12  # F17: 0 -> 12
$ws.Cells.Item(18, 6).Value = 12  # F18: 0 -> 12
$ws.Cells.Item(19, 5).Value = 1  # E19: 0 -> 1
$ws.Cells.Item(19, 6).Value = 12  # F19: 0 -> 12
$ws.Cells.Item(20, 6).Value = 12  # F20: 0 -> 12
$ws.Cells.Item(21, 6).Value = 12  # F21: 0 -> 12
$ws.Cells.Item(22, 6).Value = 12  # F22: 0 -> 12
$ws.Cells.Item(23, 6).Value = 12  # F23: 0 -> 12
$ws.Cells.Item(24, 6).Value = 2  # F24: 11 -> 2
$ws.Cells.Item(25, 5).Value = 1  # E25: 0 -> 1
$ws.Cells.Item(25, 6).Value = 2  # F25: 11 -> 2
$ws.Cells.Item(26, 6).Value = 2  # F26: 11 -> 2
$ws.Cells.Item(27, 6).Value = 2  # F27: 11 -> 2
$ws.Cells.Item(28, 6).Value = 8  # F28: 1 -> 8
$ws.Cells.Item(29, 5).Value = 1  # E29: 0 -> 1
$ws.Cells.Item(29, 6).Value = 8  # F29: 1 -> 8
$ws.Cells.Item(30, 5).Value = 1  # E30: 0 -> 1
$ws.Cells.Item(30, 6).Value = 8  # F30: 1 -> 8
$ws.Cells.Item(31, 6).Value = 8  # F31: 1 -> 8
$ws.Cells.Item(32, 6).Value = 8  # F32: 1 -> 8
$ws.Cells.Item(33, 5).Value = 1  # E33: 0 -> 1
$ws.Cells.Item(33, 6).Value = 6  # F33: 2 -> 6
$ws.Cells.Item(34, 6).Value = 6  # F34: 2 -> 6
$ws.Cells.Item(35, 6).Value = 6  # F35: 2 -> 6
$ws.Cells.Item(36, 6).Value = 13  # F36: 8 -> 13
$ws.Cells.Item(37, 5).Value = 1  # E37: 0 -> 1
$ws.Cells.Item(37, 6).Value = 13  # F37: 2 -> 13
$ws.Cells.Item(38, 5).Value = 1  # E38: 0 -> 1
$ws.Cells.Item(38, 6).Value = 13  # F38: 2 -> 13
$ws.Cells.Item(39, 6).Value = 13  # F39: 8 -> 13
$ws.Cells.Item(40, 6).Value = 6  # F40: 2 -> 6
$ws.Cells.Item(41, 6).Value = 13  # F41: 8 -> 13
$ws.Cells.Item(42, 6).Value = 13  # F42: 8 -> 13
$ws.Cells.Item(43, 5).Value = 0  # E43: 1 -> 0
$ws.Cells.Item(43, 6).Value = 13  # F43: 8 -> 13
$ws.Cells.Item(44, 5).Value = 1  # E44: 0 -> 1
$ws.Cells.Item(44, 6).Value = 13  # F44: 8 -> 13
$ws.Cells.Item(45, 6).Value = 10  # F45: 12 -> 10
$ws.Cells.Item(46, 6).Value = 10  # F46: 12 -> 10
$ws.Cells.Item(47, 6).Value = 10  # F47: 12 -> 10
$ws.Cells.Item(48, 6).Value = 10  # F48: 12 -> 10
$ws.Cells.Item(49, 5).Value = 1  # E49: 0 -> 1
$ws.Cells.Item(49, 6).Value = 10  # F49: 12 -> 10
$ws.Cells.Item(50, 6).Value = 3  # F50: 5 -> 3
$ws.Cells.Item(51, 6).Value = 10  # F51: 12 -> 10
$ws.Cells.Item(52, 6).Value = 10  # F52: 12 -> 10
$ws.Cells.Item(53, 5).Value = 1  # E53: 0 -> 1
$ws.Cells.Item(53, 6).Value = 10  # F53: 12 -> 10
$ws.Cells.Item(54, 6).Value = 10  # F54: 12 -> 10
$ws.Cells.Item(55, 6).Value = 10  # F55: 12 -> 10
$ws.Cells.Item(56, 6).Value = 1  # F56: 3 -> 1
$ws.Cells.Item(57, 6).Value = 1  # F57: 3 -> 1
$ws.Cells.Item(58, 6).Value = 1  # F58: 3 -> 1
$ws.Cells.Item(59, 5).Value = 1  # E59: 0 -> 1
$ws.Cells.Item(59, 6).Value = 1  # F59: 3 -> 1
$ws.Cells.Item(60, 5).Value = 1  # E60: 0 -> 1
$ws.Cells.Item(60, 6).Value = 1  # F60: 3 -> 1
$ws.Cells.Item(61, 6).Value = 1  # F61: 3 -> 1
$ws.Cells.Item(62, 6).Value = 1  # F62: 3 -> 1
$ws.Cells.Item(63, 6).Value = 9  # F63: 7 -> 9
$ws.Cells.Item(64, 6).Value = 9  # F64: 7 -> 9
$ws.Cells.Item(65, 6).Value = 9  # F65: 7 -> 9
$ws.Cells.Item(66, 6).Value = 9  # F66: 7 -> 9
$ws.Cells.Item(67, 6).Value = 9  # F67: 7 -> 9
$ws.Cells.Item(68, 5).Value = 1  # E68: 0 -> 1
$ws.Cells.Item(68, 6).Value = 9  # F68: 7 -> 9
$ws.Cells.Item(69, 6).Value = 9  # F69: 7 -> 9
$ws.Cells.Item(70, 5).Value = 1  # E70: 0 -> 1
$ws.Cells.Item(70, 6).Value = 9  # F70: 7 -> 9
$ws.Cells.Item(71, 5).Value = 1  # E71: 0 -> 1
$ws.Cells.Item(71, 6).Value = 11  # F71: 1 -> 11
$ws.Cells.Item(72, 6).Value = 11  # F72: 13 -> 11
$ws.Cells.Item(73, 5).Value = 1  # E73: 0 -> 1
$ws.Cells.Item(73, 6).Value = 2  # F73: 11 -> 2
$ws.Cells.Item(74, 5).Value = 1  # E74: 0 -> 1
$ws.Cells.Item(74, 6).Value = 11  # F74: 13 -> 11
$ws.Cells.Item(75, 6).Value = 11  # F75: 13 -> 11
$ws.Cells.Item(76, 5).Value = 1  # E76: 0 -> 1
$ws.Cells.Item(76, 6).Value = 11  # F76: 13 -> 11
$ws.Cells.Item(77, 6).Value = 8  # F77: 1 -> 8
$ws.Cells.Item(78, 5).Value = 0  # E78: 1 -> 0
$ws.Cells.Item(78, 6).Value = 11  # F78: 13 -> 11
$ws.Cells.Item(79, 6).Value = 0  # F79: 9 -> 0
$ws.Cells.Item(80, 6).Value = 0  # F80: 9 -> 0
$ws.Cells.Item(81, 5).Value = 1  # E81: 0 -> 1
$ws.Cells.Item(81, 6).Value = 0  # F81: 9 -> 0
$ws.Cells.Item(82, 5).Value = 1  # E82: 0 -> 1
$ws.Cells.Item(82, 6).Value = 0  # F82: 9 -> 0
$ws.Cells.Item(83, 5).Value = 1  # E83: 0 -> 1
$ws.Cells.Item(83, 6).Value = 0  # F83: 9 -> 0
$ws.Cells.Item(84, 6).Value = 0  # F84: 9 -> 0
$ws.Cells.Item(85, 6).Value = 0  # F85: 9 -> 0
$ws.Cells.Item(86, 6).Value = 0  # F86: 9 -> 0
$ws.Cells.Item(87, 6).Value = 0  # F87: 9 -> 0
$ws.Cells.Item(88, 6).Value = 7  # F88: 6 -> 7
$ws.Cells.Item(89, 6).Value = 7  # F89: 6 -> 7
$ws.Cells.Item(90, 6).Value = 7  # F90: 6 -> 7
$ws.Cells.Item(91, 6).Value = 7  # F91: 6 -> 7
$ws.Cells.Item(92, 6).Value = 7  # F92: 6 -> 7
$ws.Cells.Item(93, 5).Value = 1  # E93: 0 -> 1
$ws.Cells.Item(93, 6).Value = 7  # F93: 6 -> 7
$ws.Cells.Item(94, 6).Value = 7  # F94: 6 -> 7
$ws.Cells.Item(95, 5).Value = 1  # E95: 0 -> 1
$ws.Cells.Item(95, 6).Value = 7  # F95: 6 -> 7
$ws.Cells.Item(96, 6).Value = 5  # F96: 4 -> 5
$ws.Cells.Item(97, 5).Value = 1  # E97: 0 -> 1
$ws.Cells.Item(97, 6).Value = 5  # F97: 4 -> 5
$ws.Cells.Item(98, 6).Value = 5  # F98: 4 -> 5
$ws.Cells.Item(99, 6).Value = 5  # F99: 4 -> 5
$ws.Cells.Item(100, 6).Value = 5  # F100: 4 -> 5
$ws.Cells.Item(101, 6).Value = 5  # F101: 4 -> 5
$ws.Cells.Item(102, 6).Value = 5  # F102: 4 -> 5
$ws.Cells.Item(103, 5).Value = 1  # E103: 0 -> 1
$ws.Cells.Item(103, 6).Value = 5  # F103: 4 -> 5
$ws.Cells.Item(104, 6).Value = 6  # F104: 2 -> 6
$ws.Cells.Item(105, 5).Value = 1  # E105: 0 -> 1
$ws.Cells.Item(105, 6).Value = 6  # F105: 2 -> 6
$ws.Cells.Item(106, 6).Value = 6  # F106: 2 -> 6
$ws.Cells.Item(107, 6).Value = 13  # F107: 8 -> 13
$ws.Cells.Item(108, 6).Value = 13  # F108: 8 -> 13
$ws.Cells.Item(109, 6).Value = 6  # F109: 2 -> 6
$ws.Cells.Item(110, 6).Value = 13  # F110: 8 -> 13
$ws.Cells.Item(111, 6).Value = 6  # F111: 2 -> 6
$ws.Cells.Item(112, 6).Value = 6  # F112: 2 -> 6
$ws.Cells.Item(113, 6).Value = 6  # F113: 2 -> 6
$ws.Cells.Item(114, 6).Value = 3  # F114: 5 -> 3
$ws.Cells.Item(115, 6).Value = 3  # F115: 5 -> 3
$ws.Cells.Item(116, 6).Value = 3  # F116: 5 -> 3
$ws.Cells.Item(117, 6).Value = 3  # F117: 5 -> 3
$ws.Cells.Item(118, 6).Value = 10  # F118: 12 -> 10
$ws.Cells.Item(119, 6).Value = 3  # F119: 5 -> 3
